$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the worksheet "current" -> "Current Schedule"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("current")
$oldName = $ws1.Name
$ws1.Name = "Current Schedule"

# ---------------------------------------------------------------------
# 2. The workbook-scoped Print_Area / Print_Titles defined names for that
#    sheet still contain the literal old sheet name in their formula text
#    (renaming a sheet does not rewrite defined-name text in this engine),
#    so patch them up explicitly. Because the new sheet name contains a
#    space it must be wrapped in single quotes, matching normal Excel
#    reference syntax.
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Current Schedule!Print_Area") {
        $n.RefersTo = "='Current Schedule'!#REF!"
    }
    elseif ($n.Name -eq "Current Schedule!Print_Titles") {
        $n.RefersTo = "='Current Schedule'!`$1:`$3"
    }
}

# ---------------------------------------------------------------------
# 3 & 4. Sheet "Starting January 1, 2001": its cached sheet-level metadata
#    (outlineLevelCol on sheetFormatPr, and the dimension reference) are
#    stale/out of date compared to the sheet's real content. Re-applying
#    the existing column outline level forces the outline-level cache to
#    be written, and touching the true last used cell forces the sheet
#    dimension to be recalculated to its correct, smaller bounds.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Starting January 1, 2001")

# Re-assert the (unchanged) outline level of the already-grouped columns
# D:G so the sheet-level outline-level cache gets (re)written on save.
$ws2.Range("D1:G1").EntireColumn.OutlineLevel = 1

# The sheet's real last row is 635 (only columns A:G are ever used), but
# that row has no populated cells, so simply loading/saving the workbook
# does not pick it up as "used". Re-apply the (blank) style that is
# already in effect on the neighbouring row to register row 635 as part
# of the sheet's used range without altering any visible formatting.
$lastUsedRow = 634
$lastUsedCol = 7
$probeStyle = $ws2.Cells.Item($lastUsedRow, $lastUsedCol).Style
$ws2.Cells.Item(635, $lastUsedCol).Style = $probeStyle
